# Insert a new weekly record (row 331) into the Coliflor price sheet.
# This pushes all existing rows 331..457 down to 332..458, and fills the
# freshly-inserted row 331 with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 331 (shifts 331..457 -> 332..458)
$ws.Rows.Item(331).Insert()

# Populate the new row 331 with the new record's data
$ws.Range("A331").Value = 5
$ws.Range("B331").Value = "Macroferia Regional de Talca"
$ws.Range("C331").Value = "Maule"
$ws.Range("D331").Value = 45119
$ws.Range("E331").Value = 7
$ws.Range("F331").Value = 100112008
$ws.Range("G331").Value = "Coliflor"
$ws.Range("H331").Value = "Sin especificar"
$ws.Range("I331").Value = "Primera"
$ws.Range("J331").Value = 5000
$ws.Range("K331").Value = 600
$ws.Range("L331").Value = 600
$ws.Range("M331").Value = 600
$ws.Range("N331").Value = "`$/unidad"
$ws.Range("O331").Value = "Región del Maule"
$ws.Range("P331").Value = 600
$ws.Range("Q331").Value = 1
$ws.Range("R331").Value = "Hortaliza"
